$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit rotates the data of rows 9, 11 and 12 (row 10 is untouched):
#   new row 9  <- old row 11
#   new row 11 <- old row 12
#   new row 12 <- old row 9
# Capture the "old" values up-front (using the .Value() getter form, since
# plain property access on this host returns the reflection descriptor
# instead of invoking the getter) before any of them get overwritten.

$old9_A  = $ws.Range("A9").Value()
$old9_B  = $ws.Range("B9").Value()
$old9_D  = $ws.Range("D9").Value()
$old9_E  = $ws.Range("E9").Value()
$old9_F  = $ws.Range("F9").Value()
$old9_G  = $ws.Range("G9").Value()
$old9_H  = $ws.Range("H9").Value()
$old9_K  = $ws.Range("K9").Value()
$old9_Q  = $ws.Range("Q9").Value()
$old9_R  = $ws.Range("R9").Value()

$old11_A = $ws.Range("A11").Value()
$old11_B = $ws.Range("B11").Value()
$old11_D = $ws.Range("D11").Value()
$old11_E = $ws.Range("E11").Value()
$old11_F = $ws.Range("F11").Value()
$old11_G = $ws.Range("G11").Value()
$old11_H = $ws.Range("H11").Value()
$old11_K = $ws.Range("K11").Value()
$old11_Q = $ws.Range("Q11").Value()
$old11_R = $ws.Range("R11").Value()

$old12_A  = $ws.Range("A12").Value()
$old12_B  = $ws.Range("B12").Value()
$old12_D  = $ws.Range("D12").Value()
$old12_E  = $ws.Range("E12").Value()
$old12_F  = $ws.Range("F12").Value()
$old12_G  = $ws.Range("G12").Value()
$old12_H  = $ws.Range("H12").Value()
$old12_K  = $ws.Range("K12").Value()
$old12_Q  = $ws.Range("Q12").Value()
$old12_R  = $ws.Range("R12").Value()
$old12_AC = $ws.Range("AC12").Value()

# --- Row 9 becomes old row 11 ---
$ws.Range("A9").Value = $old11_A
$ws.Range("B9").Value = $old11_B
$ws.Range("D9").Value = $old11_D
$ws.Range("E9").Value = $old11_E
$ws.Range("F9").Value = $old11_F
$ws.Range("G9").Value = $old11_G
$ws.Range("H9").Value = $old11_H
$ws.Range("K9").Value = $old11_K
$ws.Range("Q9").Value = $old11_Q
$ws.Range("R9").Value = $old11_R

# --- Row 11 becomes old row 12 ---
$ws.Range("A11").Value = $old12_A
$ws.Range("B11").Value = $old12_B
$ws.Range("D11").Value = $old12_D
$ws.Range("E11").Value = $old12_E
$ws.Range("F11").Value = $old12_F
$ws.Range("G11").Value = $old12_G
$ws.Range("H11").Value = $old12_H
$ws.Range("K11").Value = $old12_K
$ws.Range("Q11").Value = $old12_Q
$ws.Range("R11").Value = $old12_R
# old row 12 also carried a public-comment note that row 11 didn't have yet
$ws.Range("AC11").Value = $old12_AC

# --- Row 12 becomes old row 9 ---
$ws.Range("A12").Value = $old9_A
$ws.Range("B12").Value = $old9_B
$ws.Range("D12").Value = $old9_D
$ws.Range("E12").Value = $old9_E
$ws.Range("F12").Value = $old9_F
$ws.Range("G12").Value = $old9_G
$ws.Range("H12").Value = $old9_H
$ws.Range("K12").Value = $old9_K
$ws.Range("Q12").Value = $old9_Q
$ws.Range("R12").Value = $old9_R
# row 12 no longer carries the public-comment / determination-method data
$ws.Range("AC12").ClearContents() | Out-Null
$ws.Range("AF12").ClearContents() | Out-Null
$ws.Range("J12").ClearContents() | Out-Null
$ws.Range("N12").ClearContents() | Out-Null

Write-Output "rotation applied"
